# Excel COM-interop script reproducing the authored diff:
#  - insert a new "TitlesExcelSheet" worksheet right after "loginTestData"
#  - append a new "CreateCommunity" worksheet at the end
#  - populate both new sheets with their data / shared strings
#  - move the selection on "loginTestData" and leave "CreateCommunity" as
#    the active (selected) tab, matching the saved workbook view state

$wb = $excel.ActiveWorkbook

# --- existing sheet handles -------------------------------------------------
$loginTestData = $wb.Worksheets.Item("loginTestData")

# --- 1. new sheet "TitlesExcelSheet", inserted right after loginTestData ---
$titles = $wb.Worksheets.Add($null, $loginTestData)
$titles.Name = "TitlesExcelSheet"

$titles.Range("A1").Value = "Titles"
$titles.Range("B1").Value = "CommunityName"
$titles.Range("C1").Value = "ExpectedResults"

$titles.Range("B2").Value = "argentina"
$titles.Range("C2").Value = "fail"

$titles.Range("A3").Value = "TestTeams"
$titles.Range("B3").Value = "argentina"
$titles.Range("C3").Value = "success"

$titles.Columns.Item(1).ColumnWidth = 8.27
$titles.Columns.Item(2).ColumnWidth = 13.85
$titles.Columns.Item(3).ColumnWidth = 12.58

[void]$titles.Range("C7").Select()

# --- 2. new sheet "CreateCommunity", appended at the end -------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$community = $wb.Worksheets.Add($null, $lastSheet)
$community.Name = "CreateCommunity"

$community.Range("A1").Value = "Community Name"
$community.Range("B1").Value = "ExpectedResults"

$community.Range("A2").Value = "fathiii"
$community.Range("B2").Value = "success"

$community.Columns.Item(1).ColumnWidth = 14.27
$community.Columns.Item(2).ColumnWidth = 12.58

# --- 3. update the selection left on loginTestData --------------------------
$loginTestData.Activate()
[void]$loginTestData.Range("E23").Select()

# --- 4. leave CreateCommunity as the active / selected tab ------------------
$community.Activate()
[void]$community.Range("E6").Select()

Write-Host "Sheets after edit:" ($wb.Worksheets | ForEach-Object { $_.Name })
